{"js": "// Update the LinkedIn Learning certificate list entries to append/replace the\n// completion-date suffix, per the commit's doc update.\nconst replacements = [\n  [\"Advanced Django Web Development 2016\", \"Advanced Django Web Development (Oct 2018)\"],\n  [\"Advanced Express\", \"Advanced Express (Oct 2019)\"],\n  [\"Advanced Node Js\", \"Advanced Node Js (Oct 2019)\"],\n  [\"Advanced Php Debugging Techniques\", \"Advanced Php Debugging Techniques (Jul 2018)\"],\n  [\"Advanced Python\", \"Advanced Python (Oct 2018)\"],\n  [\"Building Apis In Php Using The Slim Micro Framework\", \"Building Apis In Php Using The Slim Micro Framework (Mar 2019)\"],\n  [\"Design The Web Adding Dynamic Qr Codes\", \"Design The Web Adding Dynamic Qr Codes (Nov 2018)\"],\n  [\"Designing Restful Apis\", \"Designing Restful Apis (Oct 2018)\"],\n  [\"Extending Laravel With First Party Packages\", \"Extending Laravel With First Party Packages (Nov 2018)\"],\n  [\"Learning Django\", \"Learning Django (Oct 2018)\"],\n  [\"Learning Symfony 3\", \"Learning Symfony 3 (Nov 2018)\"],\n  [\"Node Js Essential Training 3\", \"Node Js Essential Training 3 (Oct 2019)\"],\n  [\"Node Js Microservices\", \"Node Js Microservices (Oct 2019)\"],\n  [\"Node Js Security\", \"Node Js Security (Nov 2019)\"],\n  [\"Pandas Essential Training\", \"Pandas Essential Training (Oct 2018)\"],\n  [\"Php Design Patterns\", \"Php Design Patterns (Oct 2018)\"],\n  [\"Php Testing Legacy Applications\", \"Php Testing Legacy Applications (Nov 2018)\"],\n  [\"Python Advanced Design Patterns\", \"Python Advanced Design Patterns (Nov 2018)\"],\n  [\"Typescript Essential Training\", \"Typescript Essential Training (Dec 2018)\"],\n  [\"Working Remotely 2015\", \"Working Remotely (Oct 2018)\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the LinkedIn Learning certificate list entries to append/replace the\n# completion-date suffix, per the commit's doc update.\n$d = $word.ActiveDocument\n\n$oldTexts = @(\n  \"Advanced Django Web Development 2016\",\n  \"Advanced Express\",\n  \"Advanced Node Js\",\n  \"Advanced Php Debugging Techniques\",\n  \"Advanced Python\",\n  \"Building Apis In Php Using The Slim Micro Framework\",\n  \"Design The Web Adding Dynamic Qr Codes\",\n  \"Designing Restful Apis\",\n  \"Extending Laravel With First Party Packages\",\n  \"Learning Django\",\n  \"Learning Symfony 3\",\n  \"Node Js Essential Training 3\",\n  \"Node Js Microservices\",\n  \"Node Js Security\",\n  \"Pandas Essential Training\",\n  \"Php Design Patterns\",\n  \"Php Testing Legacy Applications\",\n  \"Python Advanced Design Patterns\",\n  \"Typescript Essential Training\",\n  \"Working Remotely 2015\"\n)\n$newTexts = @(\n  \"Advanced Django Web Development (Oct 2018)\",\n  \"Advanced Express (Oct 2019)\",\n  \"Advanced Node Js (Oct 2019)\",\n  \"Advanced Php Debugging Techniques (Jul 2018)\",\n  \"Advanced Python (Oct 2018)\",\n  \"Building Apis In Php Using The Slim Micro Framework (Mar 2019)\",\n  \"Design The Web Adding Dynamic Qr Codes (Nov 2018)\",\n  \"Designing Restful Apis (Oct 2018)\",\n  \"Extending Laravel With First Party Packages (Nov 2018)\",\n  \"Learning Django (Oct 2018)\",\n  \"Learning Symfony 3 (Nov 2018)\",\n  \"Node Js Essential Training 3 (Oct 2019)\",\n  \"Node Js Microservices (Oct 2019)\",\n  \"Node Js Security (Nov 2019)\",\n  \"Pandas Essential Training (Oct 2018)\",\n  \"Php Design Patterns (Oct 2018)\",\n  \"Php Testing Legacy Applications (Nov 2018)\",\n  \"Python Advanced Design Patterns (Nov 2018)\",\n  \"Typescript Essential Training (Dec 2018)\",\n  \"Working Remotely (Oct 2018)\"\n)\n\n$matchedCount = 0\n\n# Walk every paragraph once; compare the paragraph's text (minus its\n# trailing paragraph-mark) against the old strings using a CASE-SENSITIVE\n# .Equals() check (-eq / hashtable lookups in this host are\n# case-insensitive, which would wrongly also hit the unrelated all-caps\n# \"PHP Design Patterns\" entry under Skills > Methodologies).\nforeach ($p in $d.Paragraphs) {\n  $raw = $p.Range.Text\n  $t = $raw.TrimEnd([char]13)\n  for ($i = 0; $i -lt $oldTexts.Length; $i++) {\n    if ($t.Equals($oldTexts[$i])) {\n      $p.Range.Text = $newTexts[$i]\n      $matchedCount++\n      break\n    }\n  }\n}\n\nif ($matchedCount -ne $oldTexts.Length) {\n  throw \"Expected $($oldTexts.Length) replacements, only matched $matchedCount\"\n}\n\nWrite-Output \"Replaced $matchedCount certificate entries\"\n"}
